# Updated cryptos list on Thu May 23 14:58:39 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All Price (D) cells in this sheet are stored as plain text (e.g. "7.20",
# "67.889.75" using '.' as a thousands separator), not numbers. Force text
# formatting before writing so Excel's COM layer doesn't "helpfully" coerce
# them to numeric values (which would strip trailing zeros / reformat).

# Row 2 - Bitcoin
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.889.75"
$ws.Range("E2").Value = "  -3.17%  "

# Row 3 - Ethereum
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.780.18"
$ws.Range("E3").Value = "  +0.92%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.08%  "

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "592.83"
$ws.Range("E5").Value = "  -4.25%  "

# Row 6 - Solana
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "171.01"
$ws.Range("E6").Value = "  -6.04%  "

# Row 7 - LidoStakedEther
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.779.55"
$ws.Range("E7").Value = "  +0.94%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.14%  "

# Row 9 - XRP
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.534"
$ws.Range("E9").Value = "  +0.17%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  -4.94%  "

# Row 11 - Toncoin
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.28"
$ws.Range("E11").Value = "  -1.22%  "

# Row 12 - Cardano
$ws.Range("E12").Value = "  -2.91%  "

# Row 13 - Avalanche
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "38.03"
$ws.Range("E13").Value = "  -4.89%  "

# Row 14 - ShibaInu
$ws.Range("E14").Value = "  -4.81%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.412.15"
$ws.Range("E15").Value = "  +1.07%  "

# Row 16 - WrappedEther
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.782.66"
$ws.Range("E16").Value = "  +1.23%  "

# Row 17 - WrappedBTC
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.969.48"
$ws.Range("E17").Value = "  -3.01%  "

# Row 18 - TRON
$ws.Range("E18").Value = "  -4.59%  "

# Row 19 - Polkadot
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.20"
$ws.Range("E19").Value = "  -4.92%  "

# Row 20 - Chainlink
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.93"
$ws.Range("E20").Value = "  -2.80%  "

# Row 21 - BitcoinCash
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "484.91"
$ws.Range("E21").Value = "  -3.65%  "

# Row 22 - Uniswap
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.28"
$ws.Range("E22").Value = "  -0.24%  "

# Row 23 - Polygon
$ws.Range("E23").Value = "  +1.31%  "

# Row 24 - Litecoin
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "85.55"
$ws.Range("E24").Value = "  -1.37%  "

# Row 25 - Fetch.AI
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.35"
$ws.Range("E25").Value = "  -7.90%  "

# Row 26 - PEPE
$ws.Range("E26").Value = "  +3.91%  "

# Row 27 - InternetComputer(DFINITY)
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.15"
$ws.Range("E27").Value = "  -6.35%  "

# Row 28 - RenderToken
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.10"
$ws.Range("E28").Value = "  -9.78%  "

# Row 29 - Dai
$ws.Range("E29").Value = "  -0.26%  "

# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  +0.01%  "

# Row 31 - ImmutableX
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.42"
$ws.Range("E31").Value = "  -2.05%  "

# Row 32 - EthereumClassic
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "32.42"
$ws.Range("E32").Value = "  +6.31%  "

# Row 33 - NEARProtocol
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.54"
$ws.Range("E33").Value = "  -4.54%  "

# Row 34 - Hedera
$ws.Range("E34").Value = "  -4.21%  "

# Row 35 - FirstDigitalUSD
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.998"
$ws.Range("E35").Value = "  +0.00%  "

# Row 36 - Mantle
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.996"
$ws.Range("E36").Value = "  -5.16%  "

# Row 37 - was Filecoin, now Kaspa
$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.136"
$ws.Range("E37").Value = "  -2.00%  "

# Row 38 - was Kaspa, now Filecoin
$ws.Range("B38").Value = "Filecoin"
$ws.Range("C38").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.79"
$ws.Range("E38").Value = "  -4.85%  "

# Row 39 - TheGraph
$ws.Range("E39").Value = "  -6.80%  "

# Row 40 - Bittensor
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "441.32"
$ws.Range("E40").Value = "  +3.09%  "

# Row 41 - OKB
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "48.88"
$ws.Range("E41").Value = "  -2.35%  "

# Row 42 - Stacks
$ws.Range("E42").Value = "  -4.05%  "

# Row 43 - dogwifhat
$ws.Range("E43").Value = "  -9.52%  "

# Row 44 - Cosmos
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.29"
$ws.Range("E44").Value = "  -3.56%  "

# Row 45 - Arweave
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "41.18"
$ws.Range("E45").Value = "  -7.02%  "

# Row 46 - Maker
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.840.26"
$ws.Range("E46").Value = "  -3.99%  "

# Row 48 - VeChain
$ws.Range("E48").Value = "  -3.13%  "

# Row 49 - Monero
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "137.12"
$ws.Range("E49").Value = "  +0.10%  "

# Row 50 - InjectiveProtocol
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "26.18"
$ws.Range("E50").Value = "  -3.67%  "

# Row 51 - ThetaToken
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.30"
$ws.Range("E51").Value = "  -6.75%  "
